$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 through 9 (old per-zone/export rows), now aggregated into row 2
$ws.Range("A3:B9").EntireRow.Delete()

# Update row 2 to hold the aggregated DK value
$ws.Range("A2").Value = "c_DK"
$ws.Range("B2").Value = 33539972.01128092
